# Refactor Data File Subjects Element
# Delete the "keywords" / "Data File Subjects" / "Subject Identifier" row
# (row 28) from Sheet1, shifting subsequent rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select and delete the entire row 28 (mirrors a user right-click > Delete Row)
$ws.Rows.Item(28).Delete()

# Reflect the resulting view/selection state seen after the deletion
$ws.Application.Goto($ws.Range("A7"), $false)
$ws.Rows.Item(28).Select()

$wb.Save()
